$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orlando")

# Copy the formatting (style) of the existing table (A1:C11) down to the
# new block starting at A13, so the new cells pick up the same bordered
# style (s="1") without introducing new style entries.
$ws.Range("A1:C11").Copy() | Out-Null
$ws.Range("A13:C23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 13: header row (same text as row 1) ---
$ws.Range("A13").Value = "Year"
$ws.Range("B13").Value = "Start Day"
$ws.Range("C13").Value = "End Day"

# --- Row 14: first data row (formulas un-shared, like row 2) ---
$ws.Range("A14").Value = 2014
$ws.Range("B14").Formula = '="March"&" "&ROUNDUP(14-MOD((1+A14*5/4),7),0)'
$ws.Range("C14").Formula = '="November"&" "&ROUNDUP(7-MOD((1+A14*5/4),7),0)'

# --- Rows 15-23: remaining years, with shared formulas (like rows 3-11) ---
$years = 2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022, 2023
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 15 + $i
    $ws.Cells.Item($row, 1).Value = $years[$i]
}
$ws.Range("B15:B23").Formula = '="March"&" "&ROUNDUP(14-MOD((1+A15*5/4),7),0)'
$ws.Range("C15:C23").Formula = '="November"&" "&ROUNDUP(7-MOD((1+A15*5/4),7),0)'

$ws.Range("E13").Select() | Out-Null
